$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks (and their relationships) before restructuring
$ws.Hyperlinks.Delete()

# Drop the old rows 10-23 entirely (data no longer present after re-scrape)
$ws.Rows("10:23").Delete()

# Narrow column B from 52 to 51 characters
$ws.Columns.Item(2).ColumnWidth = 50.17

# Row 2: 【SES案件多数】バックエンドエンジニア募集(Java/PHP/Python/N
$ws.Range("A2").Value = '2025-09-25 06:27:39'
$ws.Range("B2").Value = '【SES案件多数】バックエンドエンジニア募集(Java/PHP/Python/Node.js)'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5399874'
$ws.Range("G2").Value = 320
$ws.Range("H2").Value = '🔥Python ★Java ◆Node.js ○PHP'

# Row 3: 【フリーランス募集】CTビューアーソフト気道抽出機能開発
$ws.Range("A3").Value = '2025-09-25 06:27:39'
$ws.Range("B3").Value = '【フリーランス募集】CTビューアーソフト気道抽出機能開発'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5400101'
$ws.Range("G3").Value = 68
$ws.Range("H3").Value = '◆開発'

# Row 4: 【急募】SOLIDWORKS2024での機械設計と製図依頼
$ws.Range("A4").Value = '2025-09-25 06:27:39'
$ws.Range("B4").Value = '【急募】SOLIDWORKS2024での機械設計と製図依頼'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5400338'
$ws.Range("G4").Value = 25
$ws.Range("H4").ClearContents()

# Row 5: 当社CTソフトへの機能追加:気道抽出
$ws.Range("A5").Value = '2025-09-25 06:27:39'
$ws.Range("B5").Value = '当社CTソフトへの機能追加:気道抽出'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '3,000,000 円 ~ 5,000,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5400094'
$ws.Range("G5").Value = 25
$ws.Range("H5").ClearContents()

# Row 6: 【SES案件多数/リモート可】フルスタックエンジニア募集(フロント〜バック〜クラ
$ws.Range("A6").Value = '2025-09-25 06:27:39'
$ws.Range("B6").Value = '【SES案件多数/リモート可】フルスタックエンジニア募集(フロント〜バック〜クラウドまで歓迎)'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5399877'
$ws.Range("G6").Value = 25
$ws.Range("H6").ClearContents()

# Row 7: 【SES案件多数/リモート可】インフラエンジニア募集(AWS/Linux/NW設
$ws.Range("A7").Value = '2025-09-25 06:27:39'
$ws.Range("B7").Value = '【SES案件多数/リモート可】インフラエンジニア募集(AWS/Linux/NW設計・構築 等歓迎)'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5399876'
$ws.Range("G7").Value = 25
$ws.Range("H7").ClearContents()

# Row 8: 【急募】Nuxt3でのWEBページ表示速度改善依頼
$ws.Range("A8").Value = '2025-09-25 06:27:39'
$ws.Range("B8").Value = '【急募】Nuxt3でのWEBページ表示速度改善依頼'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5400231'
$ws.Range("G8").Value = 18
$ws.Range("H8").ClearContents()

# Row 9: 限定公開 PR 限定公開の仕事
$ws.Range("A9").Value = '2025-09-25 06:27:39'
$ws.Range("B9").Value = '限定公開 PR 限定公開の仕事'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5399347'
$ws.Range("G9").Value = 13
$ws.Range("H9").ClearContents()

# Re-create hyperlinks for the URL column (F2:F9), in row order so rIds line up 1-8
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5399874') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5400101') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5400338') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5400094') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5399877') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5399876') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5400231') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5399347') | Out-Null

$ws.Range("A1").Select() | Out-Null
